# Estadisticos Segundo Parcial 26 Mayo
# Updates the "2P" (second partial) grades on the Calificaciones sheet,
# the matching attendance percentages on the Asistencias sheet, and the
# recomputed failing-grade average on the Totales sheet.

$wb = $excel.ActiveWorkbook

# --- Hoja "Calificaciones": columna J (2P / F2) y columna X (Final / F2) ---
$wsCal = $wb.Worksheets.Item("Calificaciones")

$wsCal.Range("J4").Value = 7
$wsCal.Range("X4").Value = 7

$wsCal.Range("J5").Value = 8
$wsCal.Range("X5").Value = 8

$wsCal.Range("J6").Value = 10
$wsCal.Range("X6").Value = 10

$wsCal.Range("J7").Value = 10

$wsCal.Range("J8").Value = 8
$wsCal.Range("X8").Value = 8

$wsCal.Range("J9").Value = 6
$wsCal.Range("X9").Value = 8

$wsCal.Range("J10").Value = 8
$wsCal.Range("X10").Value = 7

$wsCal.Range("J11").Value = 9

$wsCal.Range("J12").Value = 10
$wsCal.Range("X12").Value = 9

$wsCal.Range("J13").Value = 9

$wsCal.Range("J14").Value = 10
$wsCal.Range("X14").Value = 10

# --- Hoja "Asistencias": columnas J y Q (2P / 3P asistencia %) ---
$wsAsi = $wb.Worksheets.Item("Asistencias")

$wsAsi.Range("J4").Value = 85.7
$wsAsi.Range("Q4").Value = 85.7

$wsAsi.Range("J5").Value = 93.90000000000001
$wsAsi.Range("Q5").Value = 93.90000000000001

$wsAsi.Range("J8").Value = 93.90000000000001
$wsAsi.Range("Q8").Value = 93.90000000000001

$wsAsi.Range("J10").Value = 89.8
$wsAsi.Range("Q10").Value = 89.8

$wsAsi.Range("J11").Value = 98
$wsAsi.Range("Q11").Value = 98

$wsAsi.Range("J14").Value = 98
$wsAsi.Range("Q14").Value = 98

# --- Hoja "Totales": promedio de reprobados ---
$wsTot = $wb.Worksheets.Item("Totales")

$wsTot.Range("H8").Value = 8.6
